$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Locate the closing paragraph ("... whole return journey.") ---
$lastPara = $d.Paragraphs.Last
$paraRange = $lastPara.Range
$paraStart = $paraRange.Start
$paraText = $paraRange.Text
$matchIdx = $paraText.IndexOf("whole return journey.")
$afterJourney = $paraStart + $matchIdx + 21   # right after "whole return journey." (before the _GoBack bookmark)

# --- Add the trailing space to "...whole return journey." (keeps xml:space="preserve") ---
$insertPoint = $d.Range($afterJourney, $afterJourney)
$insertPoint.InsertAfter(" ")
$afterSpace = $afterJourney + 1

# --- Split the paragraph right before the relocated bookmark/trailing run,
#     isolating the bookmark into its own (soon to be "Link:") paragraph ---
$splitPoint = $d.Range($afterSpace, $afterSpace)
$splitPoint.InsertParagraphBefore()

$count = $d.Paragraphs.Count
$journeyPara = $d.Paragraphs.Item($count - 1)
$bookmarkPara = $d.Paragraphs.Item($count)

# --- Rebuild the isolated paragraph as the "Link: <url>" line, bookmark moved to the end ---
$linkFrag = "<w:p $wns><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:t>Link:</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>https://streamable.com/57c99</w:t></w:r><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$bookmarkPara.Range.InsertXML($linkFrag)

# --- Insert a blank paragraph, then the bold/underlined "Video of zumo" heading,
#     both right before the "Link:" paragraph ---
$journeyPara.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Item($count)
$blankPara.Range.InsertParagraphAfter()
$videoPara = $d.Paragraphs.Item($count + 1)

$videoFrag = "<w:p $wns><w:pPr><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`">Video of </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:t>zumo</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"
$videoPara.Range.InsertXML($videoFrag)
